# Weekly data refresh: a new observation is inserted at row 451, pushing the
# existing rows 451-479 down to 452-480 (same as the source diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 451 - this shifts rows 451:479
# down to 452:480 and carries the row-451 formatting (incl. the date
# number format on column D) down onto the newly inserted row.
$ws.Rows.Item(451).Insert()

# Populate the newly inserted row 451 with the new weekly record.
$ws.Range("A451").Value = 5
$ws.Range("B451").Value = 'Macroferia Regional de Talca'
$ws.Range("C451").Value = 'Maule'
$ws.Range("D451").Value = '2023-01-05'
$ws.Range("E451").Value = 7
$ws.Range("F451").Value = 100112032
$ws.Range("G451").Value = 'Zapallo italiano'
$ws.Range("H451").Value = 'Sin especificar'
$ws.Range("I451").Value = 'Primera'
$ws.Range("J451").Value = 500
$ws.Range("K451").Value = 3000
$ws.Range("L451").Value = 3000
$ws.Range("M451").Value = 3000
$ws.Range("N451").Value = '$/caja 50 unidades'
$ws.Range("O451").Value = 'Región del Maule'
$ws.Range("P451").Value = 60
$ws.Range("Q451").Value = 50
$ws.Range("R451").Value = 'Hortaliza'
